$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Dt. Referencia" date column (G) from 45567 (2024-10-02) to
# 45568 (2024-10-03) for every data row, reflecting the new export run.
$ws.Range("G2:G274").Value = 45568

# Row 105: Saldo Previsto / Vl. Total corrected from 642.86 to 10642.86
$ws.Range("E105").Value = 10642.86
$ws.Range("H105").Value = 10642.86

# Row 109: Saldo Previsto / Vl. Total corrected from 41447.71 to 447.71
$ws.Range("E109").Value = 447.71
$ws.Range("H109").Value = 447.71

# The export file/sheet name carries the new run's timestamp.
$ws.Name = "IClientBalance-20241003-090508-"
